# Commit: "excle data updated in commit"
#
# The FlipkartInput sheet stored a real login/password pair (a phone
# number used as loginId, and a personal string used as password, the
# latter also wired up as a mailto: hyperlink). This scrubs that data:
# clears the two cells and drops the now-stale hyperlink, then leaves
# the sheet selection parked on D6 (matching the saved view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FlipkartInput")

# loginId value (was 8983844553)
$ws.Range("E2").ClearContents()

# password value (was Him@nshu77990) - also carried a mailto hyperlink
$ws.Range("F2").ClearContents()
$ws.Range("F2").Hyperlinks.Delete()

# Park the selection on D6, matching the saved sheet view.
$ws.Range("D6").Select()
